$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits at the end
#    of the "Or at least list it alphabetically." paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a brand-new list paragraph right after the
#    "Desired hex code is #96C8A2" paragraph, re-using that
#    paragraph's formatting (ListParagraph style, ilvl=1/numId=13,
#    the 1340-twip tab stop) and carrying the new note text.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Desired hex code is #96C8A2")
$anchorParagraph = $rng.Paragraphs.First
$paragraphEnd = $anchorParagraph.Range.End

$apostrophe = [char]0x2019
$newText = "Put all of the stories and emojis views into ul and li" + $apostrophe + "s like Alex did for users."

# Append a throwaway placeholder character after the new text so the
# following bookmark insertion point is NOT the very last offset of
# the paragraph (inserting a collapsed bookmark exactly at a
# paragraph-mark position gets mis-anchored) - we delete the
# placeholder again right after placing the bookmark.
$insertPoint = $d.Range($paragraphEnd, $paragraphEnd)
$insertPoint.InsertAfter("`r" + $newText + "X")

$placeholderPos = $paragraphEnd + 1 + $newText.Length
$bookmarkRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
